$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = $ws.Cells.Item(2, 4).NumberFormat()

# Row 311
$ws.Cells.Item(311, 1).Value = 11
$ws.Cells.Item(311, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(311, 3).Value = "Bíobío"
$ws.Cells.Item(311, 4).Value = 44939
$ws.Cells.Item(311, 5).Value = 8
$ws.Cells.Item(311, 6).Value = 100114013
$ws.Cells.Item(311, 7).Value = "Zanahoria"
$ws.Cells.Item(311, 8).Value = "Sin especificar"
$ws.Cells.Item(311, 9).Value = "Primera"
$ws.Cells.Item(311, 10).Value = 600
$ws.Cells.Item(311, 11).Value = 10000
$ws.Cells.Item(311, 12).Value = 11000
$ws.Cells.Item(311, 13).Value = 10500
$ws.Cells.Item(311, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(311, 15).Value = "Región de Ñuble"
$ws.Cells.Item(311, 16).Value = 525
$ws.Cells.Item(311, 17).Value = 20
$ws.Cells.Item(311, 18).Value = "Hortaliza"

# Row 312
$ws.Cells.Item(312, 1).Value = 11
$ws.Cells.Item(312, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(312, 3).Value = "Bíobío"
$ws.Cells.Item(312, 4).Value = 44939
$ws.Cells.Item(312, 5).Value = 8
$ws.Cells.Item(312, 6).Value = 100114013
$ws.Cells.Item(312, 7).Value = "Zanahoria"
$ws.Cells.Item(312, 8).Value = "Sin especificar"
$ws.Cells.Item(312, 9).Value = "Segunda"
$ws.Cells.Item(312, 10).Value = 300
$ws.Cells.Item(312, 11).Value = 9000
$ws.Cells.Item(312, 12).Value = 9000
$ws.Cells.Item(312, 13).Value = 9000
$ws.Cells.Item(312, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(312, 15).Value = "Región de Ñuble"
$ws.Cells.Item(312, 16).Value = 450
$ws.Cells.Item(312, 17).Value = 20
$ws.Cells.Item(312, 18).Value = "Hortaliza"

# Row 313
$ws.Cells.Item(313, 1).Value = 11
$ws.Cells.Item(313, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(313, 3).Value = "Bíobío"
$ws.Cells.Item(313, 4).Value = 44425
$ws.Cells.Item(313, 5).Value = 8
$ws.Cells.Item(313, 6).Value = 100114013
$ws.Cells.Item(313, 7).Value = "Zanahoria"
$ws.Cells.Item(313, 8).Value = "Sin especificar"
$ws.Cells.Item(313, 9).Value = "Primera"
$ws.Cells.Item(313, 10).Value = 600
$ws.Cells.Item(313, 11).Value = 5000
$ws.Cells.Item(313, 12).Value = 5500
$ws.Cells.Item(313, 13).Value = 5250
$ws.Cells.Item(313, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(313, 15).Value = "Región de Ñuble"
$ws.Cells.Item(313, 16).Value = 262
$ws.Cells.Item(313, 17).Value = 20
$ws.Cells.Item(313, 18).Value = "Hortaliza"

# Row 314
$ws.Cells.Item(314, 1).Value = 11
$ws.Cells.Item(314, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(314, 3).Value = "Bíobío"
$ws.Cells.Item(314, 4).Value = 44425
$ws.Cells.Item(314, 5).Value = 8
$ws.Cells.Item(314, 6).Value = 100114013
$ws.Cells.Item(314, 7).Value = "Zanahoria"
$ws.Cells.Item(314, 8).Value = "Sin especificar"
$ws.Cells.Item(314, 9).Value = "Segunda"
$ws.Cells.Item(314, 10).Value = 300
$ws.Cells.Item(314, 11).Value = 4500
$ws.Cells.Item(314, 12).Value = 4500
$ws.Cells.Item(314, 13).Value = 4500
$ws.Cells.Item(314, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(314, 15).Value = "Región de Ñuble"
$ws.Cells.Item(314, 16).Value = 225
$ws.Cells.Item(314, 17).Value = 20
$ws.Cells.Item(314, 18).Value = "Hortaliza"

# Row 315
$ws.Cells.Item(315, 1).Value = 11
$ws.Cells.Item(315, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(315, 3).Value = "Bíobío"
$ws.Cells.Item(315, 4).Value = 44827
$ws.Cells.Item(315, 5).Value = 8
$ws.Cells.Item(315, 6).Value = 100114013
$ws.Cells.Item(315, 7).Value = "Zanahoria"
$ws.Cells.Item(315, 8).Value = "Sin especificar"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 500
$ws.Cells.Item(315, 11).Value = 11000
$ws.Cells.Item(315, 12).Value = 11000
$ws.Cells.Item(315, 13).Value = 11000
$ws.Cells.Item(315, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(315, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(315, 16).Value = 550
$ws.Cells.Item(315, 17).Value = 20
$ws.Cells.Item(315, 18).Value = "Hortaliza"

# Row 316
$ws.Cells.Item(316, 1).Value = 11
$ws.Cells.Item(316, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(316, 3).Value = "Bíobío"
$ws.Cells.Item(316, 4).Value = 44827
$ws.Cells.Item(316, 5).Value = 8
$ws.Cells.Item(316, 6).Value = 100114013
$ws.Cells.Item(316, 7).Value = "Zanahoria"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Segunda"
$ws.Cells.Item(316, 10).Value = 300
$ws.Cells.Item(316, 11).Value = 9000
$ws.Cells.Item(316, 12).Value = 9000
$ws.Cells.Item(316, 13).Value = 9000
$ws.Cells.Item(316, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(316, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(316, 16).Value = 450
$ws.Cells.Item(316, 17).Value = 20
$ws.Cells.Item(316, 18).Value = "Hortaliza"

# Row 317
$ws.Cells.Item(317, 1).Value = 11
$ws.Cells.Item(317, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(317, 3).Value = "Bíobío"
$ws.Cells.Item(317, 4).Value = 44504
$ws.Cells.Item(317, 5).Value = 8
$ws.Cells.Item(317, 6).Value = 100114013
$ws.Cells.Item(317, 7).Value = "Zanahoria"
$ws.Cells.Item(317, 8).Value = "Sin especificar"
$ws.Cells.Item(317, 9).Value = "Primera"
$ws.Cells.Item(317, 10).Value = 350
$ws.Cells.Item(317, 11).Value = 6000
$ws.Cells.Item(317, 12).Value = 7000
$ws.Cells.Item(317, 13).Value = 6571
$ws.Cells.Item(317, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(317, 15).Value = "Chillán"
$ws.Cells.Item(317, 16).Value = 329
$ws.Cells.Item(317, 17).Value = 20
$ws.Cells.Item(317, 18).Value = "Hortaliza"

# Row 318
$ws.Cells.Item(318, 1).Value = 11
$ws.Cells.Item(318, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(318, 3).Value = "Bíobío"
$ws.Cells.Item(318, 4).Value = 44370
$ws.Cells.Item(318, 5).Value = 8
$ws.Cells.Item(318, 6).Value = 100114013
$ws.Cells.Item(318, 7).Value = "Zanahoria"
$ws.Cells.Item(318, 8).Value = "Sin especificar"
$ws.Cells.Item(318, 9).Value = "Primera"
$ws.Cells.Item(318, 10).Value = 600
$ws.Cells.Item(318, 11).Value = 4500
$ws.Cells.Item(318, 12).Value = 5000
$ws.Cells.Item(318, 13).Value = 4750
$ws.Cells.Item(318, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(318, 15).Value = "Región de Ñuble"
$ws.Cells.Item(318, 16).Value = 238
$ws.Cells.Item(318, 17).Value = 20
$ws.Cells.Item(318, 18).Value = "Hortaliza"

# Row 319
$ws.Cells.Item(319, 1).Value = 11
$ws.Cells.Item(319, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(319, 3).Value = "Bíobío"
$ws.Cells.Item(319, 4).Value = 44691
$ws.Cells.Item(319, 5).Value = 8
$ws.Cells.Item(319, 6).Value = 100114013
$ws.Cells.Item(319, 7).Value = "Zanahoria"
$ws.Cells.Item(319, 8).Value = "Sin especificar"
$ws.Cells.Item(319, 9).Value = "Primera"
$ws.Cells.Item(319, 10).Value = 600
$ws.Cells.Item(319, 11).Value = 7000
$ws.Cells.Item(319, 12).Value = 7500
$ws.Cells.Item(319, 13).Value = 7250
$ws.Cells.Item(319, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(319, 15).Value = "Región de Ñuble"
$ws.Cells.Item(319, 16).Value = 362
$ws.Cells.Item(319, 17).Value = 20
$ws.Cells.Item(319, 18).Value = "Hortaliza"

# Row 320
$ws.Cells.Item(320, 1).Value = 11
$ws.Cells.Item(320, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(320, 3).Value = "Bíobío"
$ws.Cells.Item(320, 4).Value = 44691
$ws.Cells.Item(320, 5).Value = 8
$ws.Cells.Item(320, 6).Value = 100114013
$ws.Cells.Item(320, 7).Value = "Zanahoria"
$ws.Cells.Item(320, 8).Value = "Sin especificar"
$ws.Cells.Item(320, 9).Value = "Segunda"
$ws.Cells.Item(320, 10).Value = 300
$ws.Cells.Item(320, 11).Value = 6000
$ws.Cells.Item(320, 12).Value = 6000
$ws.Cells.Item(320, 13).Value = 6000
$ws.Cells.Item(320, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(320, 15).Value = "Región de Ñuble"
$ws.Cells.Item(320, 16).Value = 300
$ws.Cells.Item(320, 17).Value = 20
$ws.Cells.Item(320, 18).Value = "Hortaliza"

# Row 321
$ws.Cells.Item(321, 1).Value = 11
$ws.Cells.Item(321, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(321, 3).Value = "Bíobío"
$ws.Cells.Item(321, 4).Value = 44306
$ws.Cells.Item(321, 5).Value = 8
$ws.Cells.Item(321, 6).Value = 100114013
$ws.Cells.Item(321, 7).Value = "Zanahoria"
$ws.Cells.Item(321, 8).Value = "Sin especificar"
$ws.Cells.Item(321, 9).Value = "Primera"
$ws.Cells.Item(321, 10).Value = 400
$ws.Cells.Item(321, 11).Value = 6500
$ws.Cells.Item(321, 12).Value = 6500
$ws.Cells.Item(321, 13).Value = 6500
$ws.Cells.Item(321, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(321, 15).Value = "Región de Ñuble"
$ws.Cells.Item(321, 16).Value = 325
$ws.Cells.Item(321, 17).Value = 20
$ws.Cells.Item(321, 18).Value = "Hortaliza"

# Row 322
$ws.Cells.Item(322, 1).Value = 11
$ws.Cells.Item(322, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(322, 3).Value = "Bíobío"
$ws.Cells.Item(322, 4).Value = 44306
$ws.Cells.Item(322, 5).Value = 8
$ws.Cells.Item(322, 6).Value = 100114013
$ws.Cells.Item(322, 7).Value = "Zanahoria"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Segunda"
$ws.Cells.Item(322, 10).Value = 400
$ws.Cells.Item(322, 11).Value = 5000
$ws.Cells.Item(322, 12).Value = 5000
$ws.Cells.Item(322, 13).Value = 5000
$ws.Cells.Item(322, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(322, 15).Value = "Región de Ñuble"
$ws.Cells.Item(322, 16).Value = 250
$ws.Cells.Item(322, 17).Value = 20
$ws.Cells.Item(322, 18).Value = "Hortaliza"

# Row 323
$ws.Cells.Item(323, 1).Value = 11
$ws.Cells.Item(323, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(323, 3).Value = "Bíobío"
$ws.Cells.Item(323, 4).Value = 44812
$ws.Cells.Item(323, 5).Value = 8
$ws.Cells.Item(323, 6).Value = 100114013
$ws.Cells.Item(323, 7).Value = "Zanahoria"
$ws.Cells.Item(323, 8).Value = "Sin especificar"
$ws.Cells.Item(323, 9).Value = "Primera"
$ws.Cells.Item(323, 10).Value = 800
$ws.Cells.Item(323, 11).Value = 9000
$ws.Cells.Item(323, 12).Value = 10000
$ws.Cells.Item(323, 13).Value = 9500
$ws.Cells.Item(323, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(323, 15).Value = "Región de Ñuble"
$ws.Cells.Item(323, 16).Value = 475
$ws.Cells.Item(323, 17).Value = 20
$ws.Cells.Item(323, 18).Value = "Hortaliza"

# Row 324
$ws.Cells.Item(324, 1).Value = 11
$ws.Cells.Item(324, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(324, 3).Value = "Bíobío"
$ws.Cells.Item(324, 4).Value = 44812
$ws.Cells.Item(324, 5).Value = 8
$ws.Cells.Item(324, 6).Value = 100114013
$ws.Cells.Item(324, 7).Value = "Zanahoria"
$ws.Cells.Item(324, 8).Value = "Sin especificar"
$ws.Cells.Item(324, 9).Value = "Segunda"
$ws.Cells.Item(324, 10).Value = 400
$ws.Cells.Item(324, 11).Value = 8000
$ws.Cells.Item(324, 12).Value = 8000
$ws.Cells.Item(324, 13).Value = 8000
$ws.Cells.Item(324, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(324, 15).Value = "Región de Ñuble"
$ws.Cells.Item(324, 16).Value = 400
$ws.Cells.Item(324, 17).Value = 20
$ws.Cells.Item(324, 18).Value = "Hortaliza"
$ws.Cells.Item(324, 4).NumberFormat = $dateFormat

# Row 325
$ws.Cells.Item(325, 1).Value = 11
$ws.Cells.Item(325, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(325, 3).Value = "Bíobío"
$ws.Cells.Item(325, 4).Value = 44791
$ws.Cells.Item(325, 5).Value = 8
$ws.Cells.Item(325, 6).Value = 100114013
$ws.Cells.Item(325, 7).Value = "Zanahoria"
$ws.Cells.Item(325, 8).Value = "Sin especificar"
$ws.Cells.Item(325, 9).Value = "Primera"
$ws.Cells.Item(325, 10).Value = 250
$ws.Cells.Item(325, 11).Value = 14000
$ws.Cells.Item(325, 12).Value = 15000
$ws.Cells.Item(325, 13).Value = 14600
$ws.Cells.Item(325, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(325, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(325, 16).Value = 730
$ws.Cells.Item(325, 17).Value = 20
$ws.Cells.Item(325, 18).Value = "Hortaliza"
$ws.Cells.Item(325, 4).NumberFormat = $dateFormat
